# Fix #12728 - [Bug] Update ExportCiteo.xlsx
$wb = $excel.ActiveWorkbook

$wsEco = $wb.Worksheets.Item("Eco emballage")
$wsMat = $wb.Worksheets.Item("Materials")

# ---------------------------------------------------------------------------
# 1. "Eco emballage" sheet: fix the VLOOKUP helper formulas in row 2 (K2:AA2)
#    so they reference the (now English) "Materials" sheet instead of the
#    stale French "Matériaux" sheet name.
# ---------------------------------------------------------------------------
$cols = @("K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")
foreach ($c in $cols) {
    $cell = $c + "2"
    $text = 'excel|IFERROR(IF(G2="","",VLOOKUP(E2&"C"&INDIRECT("' + $c + '"&ROW($A$2)+5)&"P",Materials!$F$4:$G$10000,2,FALSE)),"")'
    $wsEco.Range($cell).Value = $text
}

# ---------------------------------------------------------------------------
# 2. "Eco emballage" sheet: relabel several packaging category names
#    (row 5 and row 6 legend cells).
# ---------------------------------------------------------------------------
$wsEco.Range("S5").Value = "Rigid"
$wsEco.Range("V5").Value = "Flexible"

$wsEco.Range("M6").Value = "Paper -cardboard other than bricks"
$wsEco.Range("P6").Value = "Clear PET Bottle"
$wsEco.Range("Q6").Value = "Dark/colored PET, PE, PP Bottle"
$wsEco.Range("S6").Value = "Rigid packaging PE, PP or PET"
$wsEco.Range("T6").Value = "Rigid PS packaging"
$wsEco.Range("V6").Value = "Flexible PE packaging"
$wsEco.Range("X6").Value = "Packaging containing PVC"

# ---------------------------------------------------------------------------
# 3. "Materials" sheet: autofilter only covers the header row (B3:H3) instead
#    of spuriously extending one row further than the sheet's data (B3:H4).
# ---------------------------------------------------------------------------
$wsMat.Range("B3:H4").AutoFilter() | Out-Null
$wsMat.Range("B3:H3").AutoFilter() | Out-Null

# ---------------------------------------------------------------------------
# 3b. Workbook-level "_FilterDatabase" names: the hidden/stale entry must
#     point at $H$3 and the live (visible) one at $H$4, matching the
#     autofilter range restored above.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        if ($n.Visible) {
            $n.RefersTo = "=Materials!`$B`$3:`$H`$4"
        } else {
            $n.RefersTo = "=Materials!`$B`$3:`$H`$3"
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Restore the previously-active selections on both sheets, matching the
#    state captured in the fixed workbook ("Eco emballage" ends up the
#    active/selected sheet).
# ---------------------------------------------------------------------------
$wsMat.Activate()
$wsMat.Range("U21").Select() | Out-Null

$wsEco.Activate()
$wsEco.Range("S6").Select() | Out-Null

$wb.Save()
